# Consolidate fragmented text runs (one run per word) into a single run
# per paragraph, for the slide titles and the "(An/an) image" captions.
#
# Since the final rendered text is unchanged, we first set a throwaway
# value so the underlying text actually changes, then set the real
# target text -- this forces the writer to re-emit the paragraph as a
# single consolidated <a:r> run instead of leaving the original
# (already-matching) fragmented runs untouched.

$p = $ppt.ActivePresentation

function Set-ConsolidatedText($shape, [string]$text) {
    $shape.TextFrame.TextRange.Text = "."
    $shape.TextFrame.TextRange.Text = $text
}

Set-ConsolidatedText $p.Slides.Item(1).Shapes.Item(1) "Slide 1 (Content)"
Set-ConsolidatedText $p.Slides.Item(2).Shapes.Item(1) "Slide 2 (Content)"
Set-ConsolidatedText $p.Slides.Item(3).Shapes.Item(1) "Slide 3 (Content)"
Set-ConsolidatedText $p.Slides.Item(4).Shapes.Item(1) "Slide 4 (Content)"
Set-ConsolidatedText $p.Slides.Item(5).Shapes.Item(1) "Slide 5 (Two Content)"

Set-ConsolidatedText $p.Slides.Item(6).Shapes.Item(1) "Slide 6 (Two Content Right)"
Set-ConsolidatedText $p.Slides.Item(6).Shapes.Item(3) "an image"

Set-ConsolidatedText $p.Slides.Item(7).Shapes.Item(1) "Slide 7 (Content with Caption)"
Set-ConsolidatedText $p.Slides.Item(7).Shapes.Item(4) "An image"

Set-ConsolidatedText $p.Slides.Item(8).Shapes.Item(1) "Slide 8 (Comparison)"
Set-ConsolidatedText $p.Slides.Item(8).Shapes.Item(4) "An image"

Set-ConsolidatedText $p.Slides.Item(9).Shapes.Item(1) "Slide 10 (Content)"
Set-ConsolidatedText $p.Slides.Item(10).Shapes.Item(1) "Slide 11 (Content)"
Set-ConsolidatedText $p.Slides.Item(11).Shapes.Item(1) "Slide 12 (Content)"
